$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.166.93'
$ws.Range("E2").Value = '  +3.41%  '
$ws.Range("D3").Value = '3.115.37'
$ws.Range("E3").Value = '  +1.43%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Formula = "'523.59"
$ws.Range("E5").Value = '  +1.76%  '
$ws.Range("D6").Formula = "'145.06"
$ws.Range("E6").Value = '  +2.78%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +1.43%  '
$ws.Range("E9").Value = '  +3.01%  '
$ws.Range("E10").Value = '  +1.37%  '
$ws.Range("E11").Value = '  +3.57%  '
$ws.Range("D12").Value = '3.652.92'
$ws.Range("E12").Value = '  +1.54%  '
$ws.Range("E13").Value = '  +1.64%  '
$ws.Range("D14").Formula = "'27.36"
$ws.Range("E14").Value = '  +7.40%  '
$ws.Range("E15").Value = '  +2.17%  '
$ws.Range("D16").Value = '59.124.46'
$ws.Range("E16").Value = '  +3.25%  '
$ws.Range("E17").Value = '  +5.57%  '
$ws.Range("D18").Value = '3.115.79'
$ws.Range("E18").Value = '  +1.26%  '
$ws.Range("D19").Formula = "'13.10"
$ws.Range("E19").Value = '  +0.61%  '
$ws.Range("D20").Formula = "'8.31"
$ws.Range("E20").Value = '  +1.93%  '
$ws.Range("D21").Formula = "'340.95"
$ws.Range("E21").Value = '  +1.39%  '
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("E23").Value = '  +2.14%  '
$ws.Range("D24").Formula = "'66.07"
$ws.Range("E24").Value = '  +1.27%  '
$ws.Range("E25").Value = '  +2.38%  '
$ws.Range("D26").Formula = "'0.997"
$ws.Range("E26").Value = '  -0.47%  '
$ws.Range("E27").Value = '  -1.36%  '
$ws.Range("D28").Formula = "'6.70"
$ws.Range("E28").Value = '  +3.62%  '
$ws.Range("E29").Value = '  +4.18%  '
$ws.Range("E30").Value = '  +2.38%  '
$ws.Range("E31").Value = '  +3.84%  '
$ws.Range("D32").Formula = "'21.14"
$ws.Range("E32").Value = '  +2.06%  '
$ws.Range("D33").Formula = "'155.64"
$ws.Range("E33").Value = '  +1.07%  '
$ws.Range("D34").Formula = "'4.69"
$ws.Range("E34").Value = '  +3.21%  '
$ws.Range("D35").Formula = "'6.19"
$ws.Range("E35").Value = '  +5.99%  '
$ws.Range("D36").Formula = "'27.38"
$ws.Range("E36").Value = '  +4.41%  '
$ws.Range("E37").Value = '  +6.23%  '
$ws.Range("E38").Value = '  +2.41%  '
$ws.Range("E39").Value = '  +3.17%  '
$ws.Range("D40").Value = '3.160.00'
$ws.Range("E40").Value = '  +1.47%  '
$ws.Range("D41").Formula = "'36.93"
$ws.Range("E41").Value = '  -0.19%  '
$ws.Range("D42").Formula = "'1.00"
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("E43").Value = '  -0.58%  '
$ws.Range("D44").Formula = "'1.47"
$ws.Range("E44").Value = '  +6.47%  '
$ws.Range("D45").Value = '2.294.13'
$ws.Range("D46").Formula = "'0.0260"
$ws.Range("E46").Value = '  +2.62%  '
$ws.Range("D47").Formula = "'21.14"
$ws.Range("E47").Value = '  +5.61%  '
$ws.Range("D48").Formula = "'0.965"
$ws.Range("E48").Value = '  +1.78%  '
$ws.Range("D49").Formula = "'6.03"
$ws.Range("E49").Value = '  +3.02%  '
$ws.Range("D50").Formula = "'0.759"
$ws.Range("E50").Value = '  +10.55%  '
$ws.Range("D51").Formula = "'262.14"
$ws.Range("E51").Value = '  +11.22%  '

Write-Output "Updated 79 cells"
